$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after Sheet1
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Row 1
$ws2.Range("A1").Value = "d65as4"
$ws2.Range("B1").Value = "das4d6"
$ws2.Range("C1").Value = "dasdsa"

# Row 3 (row 2 intentionally left empty)
$ws2.Range("A3").Value = "6a6dsa5hfghfg"
$ws2.Range("B3").Value = "tertew"
$ws2.Range("C3").Value = "gdfg"

# Match the selection / active state recorded in the edited workbook
$ws2.Range("A4").Select()
$ws2.Activate()
